# Update countries & provincias Spain
# Applies the "6 de Mayo de 2020" data refresh:
#  - timestamp label text updated (21:03 -> 21:33)
#  - a handful of country rows got refreshed totals
#  - Congo/Ruanda swapped rank (row 134/135), Gibraltar/Zambia swapped rank (row 147/148)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 21:33"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1249791
$ws.Range("C4").Value = 12158
$ws.Range("D4").Value = 204872
$ws.Range("E4").Value = 971286
$ws.Range("F4").Value = 15843
$ws.Range("G4").Value = 1362
$ws.Range("H4").Value = 73633

# --- Row 113: Maldivas ---
$ws.Range("B113").Value = 617
$ws.Range("C113").Value = 44
$ws.Range("H113").Value = 1

# --- Row 128: Estado de Palestina ---
$ws.Range("B128").Value = 374
$ws.Range("C128").Value = 3
$ws.Range("E128").Value = 198

# --- Rows 134/135: Congo and Ruanda swap rank (Ruanda now above Congo) ---
$ws.Range("A134").Value = "Ruanda"
$ws.Range("B134").Value = 268
$ws.Range("C134").Value = 7
$ws.Range("D134").Value = 130
$ws.Range("E134").Value = 138
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

$ws.Range("A135").Value = "Congo"
$ws.Range("B135").Value = 264
$ws.Range("C135").Value = 28
$ws.Range("D135").Value = 30
$ws.Range("E135").Value = 224
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 10

# --- Rows 147/148: Gibraltar and Zambia swap rank (Zambia now above Gibraltar) ---
$ws.Range("A147").Value = "Zambia"
$ws.Range("B147").Value = 146
$ws.Range("C147").Value = 8
$ws.Range("D147").Value = 101
$ws.Range("E147").Value = 41
$ws.Range("F147").Value = 1
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 4

$ws.Range("A148").Value = "Gibraltar"
$ws.Range("B148").Value = 144
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 136
$ws.Range("E148").Value = 8
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 0
